$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header in I1
$ws.Range("I1").Value = "score"

# Add the new second row of data
$ws.Range("B2").Value = "yes"
$ws.Range("C2").Value = 119
$ws.Range("D2").Value = "full body"

# Rename the F1 header from "rem_sleep" to "awakenings"
$ws.Range("F1").Value = "awakenings"

$ws.Range("E11").Select()
$excel.ActiveWindow.Zoom = 179
